$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fixed a bug while editing the related invoice: G2's total_amount was
# wrong (1208 -> 0) ---
$ws.Range("G2").Value = 0

# --- Set the "best fit" column widths (A:J) that Excel computed for this
# sheet's content (matches the <cols> block added to the worksheet XML) ---
$ws.Columns.Item(1).ColumnWidth  = 9.666666666666666   # A -> 10.5703125
$ws.Columns.Item(2).ColumnWidth  = 11.833333333333334  # B -> 12.7109375
$ws.Columns.Item(3).ColumnWidth  = 7.166666666666667   # C -> 8
$ws.Columns.Item(4).ColumnWidth  = 11.166666666666666  # D -> 12
$ws.Columns.Item(5).ColumnWidth  = 11.166666666666666  # E -> 12
$ws.Columns.Item(6).ColumnWidth  = 18.5                # F -> 19.28515625
$ws.Columns.Item(7).ColumnWidth  = 12.333333333333334  # G -> 13.140625
$ws.Columns.Item(8).ColumnWidth  = 6.5                 # H -> 7.28515625
$ws.Columns.Item(9).ColumnWidth  = 5.666666666666667   # I -> 6.42578125
$ws.Columns.Item(10).ColumnWidth = 9.5                 # J -> 10.28515625

# --- Move the active selection from G7 to H9 ---
$ws.Range("H9").Select() | Out-Null
